$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replyText = "Geachte afzender,`nDank u voor uw bericht. Kunt u meer details geven over wat u precies wilt dat we oppakken? Zo kunnen we u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam van het bedrijf]"

$ws.Range("A8").Value = "Testmail #1: Zou jij dit even op kunnen pakken?"
$ws.Range("B8").Value = $replyText
$ws.Range("C8").Value = "Zou jij dit even op kunnen pakken?"
$ws.Range("D8").Value = "mailmind.test@zohomail.eu"
$ws.Range("E8").Value = "Overig"
$ws.Range("F8").Value = "2025-08-05 19:21:21"
$ws.Range("G8").Value = "Ja"
$ws.Range("H8").Value = "Nee"
$ws.Range("I8").Value = "Ja"
$ws.Range("J8").Value = "Nee"

$ws.Rows.Item(8).AutoFit()
